$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..43 (Timestamp, IGCC Import (MW), IGCC Export (MW))
$data = @(
    @(45889, 0.242, 18.302),
    @(45889.01041666666, 0, 28.622),
    @(45889.02083333334, 0, 36.249),
    @(45889.03125, 0, 24.025),
    @(45889.04166666666, 0, 55.354),
    @(45889.05208333334, 0, 60.823),
    @(45889.0625, 0, 70.39700000000001),
    @(45889.07291666666, 0, 52.023),
    @(45889.08333333334, 0, 21.686),
    @(45889.09375, 0, 12.902),
    @(45889.10416666666, 0, 10.701),
    @(45889.11458333334, 0.02, 14.472),
    @(45889.125, 2.201, 0.386),
    @(45889.13541666666, 5.534, 0.019),
    @(45889.14583333334, 13.691, 0),
    @(45889.15625, 19.063, 0),
    @(45889.16666666666, 14.3, 0),
    @(45889.17708333334, 2.424, 4.863),
    @(45889.1875, 0.305, 3.195),
    @(45889.19791666666, 0.164, 2.165),
    @(45889.20833333334, 0.049, 3.517),
    @(45889.21875, 0.352, 2.29),
    @(45889.22916666666, 0.061, 5.013),
    @(45889.23958333334, 0.005, 10.355),
    @(45889.25, 0, 20.577),
    @(45889.26041666666, 0, 24.493),
    @(45889.27083333334, 0, 22.957),
    @(45889.28125, 0, 37.779),
    @(45889.29166666666, 0, 11.97),
    @(45889.30208333334, 0, 19.975),
    @(45889.3125, 0, 11.09),
    @(45889.32291666666, 0, 29.869),
    @(45889.33333333334, 0.004, 14.529),
    @(45889.34375, 0.021, 8.965999999999999),
    @(45889.35416666666, 0.5669999999999999, 8.423999999999999),
    @(45889.36458333334, 0, 38.742),
    @(45889.375, 0, 23.279),
    @(45889.38541666666, 0.223, 6.356),
    @(45889.39583333334, 1.574, 1.143),
    @(45889.40625, 0, 6.026),
    @(45889.41666666666, 0, 21.595),
    @(45889.42708333334, 0, 27.518)
)

$firstDataRow = 2
$lastExistingRow = 29
$lastNewRow = $firstDataRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $firstDataRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# New rows need the same date/time number format as column A of existing rows
$newARange = $ws.Range("A$($lastExistingRow + 1):A$lastNewRow")
$newARange.NumberFormat = "YYYY-MM-DD HH:MM:SS"
